# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
#
# Columns D (Price) and E (Volume(1h)) hold plain text, not numbers (some
# prices use "." as a thousands separator, and the volume cells keep
# surrounding padding spaces). Several of the new price strings (e.g.
# "568.61", "1.00") are themselves valid numeric literals, so a plain
# `.Value = "..."` assignment would make Excel auto-coerce them into real
# numbers. To keep them as text -- matching the original cell type -- those
# cells are briefly switched to the Text number format, written, then
# switched back to the default "Normal" style so no visible formatting
# change is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.222.29"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "2.550.77"
$ws.Range("E3").Value = "  +3.71%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.78%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").Value = "2.549.14"
$ws.Range("E9").Value = "  +3.71%  "

$ws.Range("E10").Value = "  +0.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.91%  "

$ws.Range("D15").Value = "3.005.79"
$ws.Range("E15").Value = "  +3.70%  "

$ws.Range("D16").Value = "63.133.47"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("E17").Value = "  +2.39%  "

$ws.Range("D18").Value = "2.525.24"
$ws.Range("E18").Value = "  +2.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.54%  "

$ws.Range("E26").Value = "  -1.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.80%  "

$ws.Range("D31").Value = "0.0₃0822"
$ws.Range("E31").Value = "  +3.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("E34").Value = "  +4.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "412.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.398"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.33%  "

$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.30%  "

$ws.Range("E42").Value = "  -3.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "153.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.605"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0962"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0524"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.83%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0239"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.90%  "
